# Update "想去人数" (F) and "最低票价" (G) values on the "展览" and "全部类型"
# sheets to reflect the latest scraped numbers.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 266
$ws1.Range("F3").Value = 71
$ws1.Range("F5").Value = 6658
$ws1.Range("F6").Value = 5435
$ws1.Range("F8").Value = 68
$ws1.Range("G9").Value = 29.9
$ws1.Range("F10").Value = 65
$ws1.Range("F11").Value = 237
$ws1.Range("F12").Value = 79

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 266
$ws4.Range("F3").Value = 71
$ws4.Range("F5").Value = 6658
$ws4.Range("F6").Value = 5435
$ws4.Range("F8").Value = 68
$ws4.Range("G9").Value = 29.9
$ws4.Range("F10").Value = 65
$ws4.Range("F11").Value = 237
$ws4.Range("F14").Value = 80
